$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GoalPoseX (row 12) and GoalPoseY (row 13) values
$ws.Range("B12").Value = 154
$ws.Range("B13").Value = -339

# Move the active selection to F31 to match the final cursor position
$ws.Range("F31").Select() | Out-Null
